$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '''245.60'
$ws.Range("E2").Formula = '''-0.49%'
$ws.Range("D3").Formula = '''27.14'
$ws.Range("E3").Formula = '''3.15%'
$ws.Range("D4").Formula = '''5.107'
$ws.Range("E4").Formula = '''0.69%'
$ws.Range("D5").Formula = '''0.05702'
$ws.Range("E5").Formula = '''1.75%'
$ws.Range("D6").Formula = '''6.512'
$ws.Range("E6").Formula = '''0.56%'
$ws.Range("E7").Formula = '''0.71%'
$ws.Range("D8").Formula = '''0.8609'
$ws.Range("E8").Formula = '''2.00%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Formula = '''0.1333'
$ws.Range("E9").Formula = '''-0.35%'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Formula = '''0.06945'
$ws.Range("E10").Formula = '''-0.56%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Formula = '''0.02844'
$ws.Range("E11").Formula = '''0.40%'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Formula = '''0.09389'
$ws.Range("E12").Formula = '''-0.29%'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Formula = '''0.001515'
$ws.Range("E13").Formula = '''-0.57%'
$ws.Range("B14").Value = 'CoinExToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D14").Formula = '''0.04037'
$ws.Range("E14").Formula = '''-13.22%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Formula = '''0.01007'
$ws.Range("E15").Formula = '''1,580.94%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Formula = '''0.006210'
$ws.Range("E16").Formula = '''0.87%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Formula = '''3.511'
$ws.Range("E17").Formula = '''-2.63%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Formula = '''3.008'
$ws.Range("E18").Formula = '''-0.13%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Formula = '''2.316'
$ws.Range("E19").Formula = '''12.67%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Formula = '''0.3165'
$ws.Range("E20").Formula = '''1.25%'
$ws.Range("E21").Formula = '''1.06%'
$ws.Range("E22").Formula = '''-1.77%'
$ws.Range("D23").Formula = '''3.583'
$ws.Range("E23").Formula = '''-4.49%'
$ws.Range("E24").Formula = '''1.73%'
$ws.Range("D25").Formula = '''0.001214'
$ws.Range("E25").Formula = '''-2.77%'
$ws.Range("D26").Formula = '''0.004476'
$ws.Range("E26").Formula = '''-2.41%'
$ws.Range("D27").Formula = '''0.00009891'
$ws.Range("E27").Formula = '''3.05%'
$ws.Range("E28").Formula = '''-25.29%'
$ws.Range("D40").Formula = '''0.03733'
$ws.Range("E40").Formula = '''1.83%'
$ws.Range("D41").Formula = '''0.006015'
$ws.Range("E41").Formula = '''-2.50%'
$ws.Range("D42").Formula = '''0.1058'
$ws.Range("E42").Formula = '''-0.17%'
$ws.Range("D43").Formula = '''0.002521'
$ws.Range("E43").Formula = '''0.84%'
$ws.Range("D44").Formula = '''0.009714'
$ws.Range("E44").Formula = '''17.64%'
$ws.Range("D45").Formula = '''0.00005150'
$ws.Range("E45").Formula = '''-4.42%'
$ws.Range("E46").Formula = '''-0.07%'
$ws.Range("D48").Formula = '''0.002506'
$ws.Range("E48").Formula = '''-3.55%'
$ws.Range("E49").Formula = '''-0.07%'
$ws.Range("E50").Formula = '''-0.07%'
